$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @('E', 'F', 'G', 'H', 'I', 'J', 'K', 'L', 'M', 'N', 'O', 'P', 'Q', 'R', 'S', 'T')

$rows = @(
    @{ Row = 2; D = "ECs"; Vals = @(3, 1, 0.4150643333333333, 1.245193, 0.02396302145531912, 0.02396302145531911, 2, 1, 10.836393, 21.672786, 0.01846862257356514, 0.01262296693390161, 4.497800236283, 26.986801417698, 0.0004425639989805323, 0.0003024844274668679) },
    @{ Row = 3; D = "FAPs"; Vals = @(3, 1, 0.4150643333333333, 1.245193, 0.02396302145531912, 0.02396302145531911, 3, 1, 139.6948166666666, 419.0844499999999, 0.2380839126543345, 0.2440890227431923, 57.98233594987221, 521.8410235488499, 0.005705209907102141, 0.005849110489002993) },
    @{ Row = 4; D = "Inflammatory-Mac"; Vals = @(3, 1, 0.4150643333333333, 1.245193, 0.02396302145531912, 0.02396302145531911, 3, 1, 150.12088, 450.36264, 0.2558532043948076, 0.2623065033256284, 62.30982297661333, 560.78840678952, 0.006131015826324923, 0.006285656367061769) },
    @{ Row = 5; D = "MuSCs"; Vals = @(3, 1, 0.4150643333333333, 1.245193, 0.02396302145531912, 0.02396302145531911, 2, 1, 32.469223, 64.938446, 0.05533777012737728, 0.03782235733776705, 13.47681639834633, 80.86089839007799, 0.001326060172851859, 0.0009063379603756581) },
    @{ Row = 6; D = "Neutrophils"; Vals = @(3, 1, 0.4150643333333333, 1.245193, 0.02396302145531912, 0.02396302145531911, 3, 1, 52.09024033333333, 156.270721, 0.08877815602319267, 0.09101737745760805, 21.62080087712811, 194.587207894153, 0.002127392857547434, 0.00218105136882354) },
    @{ Row = 7; D = $null; Vals = @(3, 1, 0.4150643333333333, 1.245193, 0.02396302145531912, 0.02396302145531911, 3, 1, 201.5345866666667, 604.60376, 0.3434783342267227, 0.3521417722019025, 83.64981885840888, 752.84836972568, 0.008230778692512227, 0.008438380842588286) },
    @{ Row = 8; D = "ECs"; Vals = @(3, 1, 16.787094, 50.361282, 0.9691738397849783, 0.9691738397849782, 2, 1, 10.836393, 21.672786, 0.01846862257356514, 0.01262296693390161, 181.911547911942, 1091.469287471652, 0.01789930585516165, 0.01223384933280823) },
    @{ Row = 9; D = "FAPs"; Vals = @(3, 1, 16.787094, 50.361282, 0.9691738397849783, 0.9691738397849782, 3, 1, 139.6948166666666, 419.0844499999999, 0.2380839126543345, 0.2440890227431923, 2345.070018696099, 21105.6301682649, 0.2307446998182328, 0.2365646954213826) },
    @{ Row = 10; D = "Inflammatory-Mac"; Vals = @(3, 1, 16.787094, 50.361282, 0.9691738397849783, 0.9691738397849782, 3, 1, 150.12088, 450.36264, 0.2558532043948076, 0.2623065033256284, 2520.09332392272, 22680.83991530448, 0.2479662325246066, 0.2542206010286704) },
    @{ Row = 11; D = "MuSCs"; Vals = @(3, 1, 16.787094, 50.361282, 0.9691738397849783, 0.9691738397849782, 2, 1, 32.469223, 64.938446, 0.05533777012737728, 0.03782235733776705, 545.063898607962, 3270.383391647772, 0.0536319191594887, 0.03665643929076324) },
    @{ Row = 12; D = "Neutrophils"; Vals = @(3, 1, 16.787094, 50.361282, 0.9691738397849783, 0.9691738397849782, 3, 1, 52.09024033333333, 156.270721, 0.08877815602319267, 0.09101737745760805, 874.4437609582578, 7869.993848624322, 0.08604146636202753, 0.0882116611977487) },
    @{ Row = 13; D = $null; Vals = @(3, 1, 16.787094, 50.361282, 0.9691738397849783, 0.9691738397849782, 3, 1, 201.5345866666667, 604.60376, 0.3434783342267227, 0.3521417722019025, 3383.18005062448, 30448.62045562032, 0.332890216065461, 0.341286593513605) },
    @{ Row = 14; D = "ECs"; Vals = @(2, 0.6666666666666666, 0.1188766666666667, 0.35663, 0.006863138759702679, 0.006863138759702677, 2, 1, 10.836393, 21.672786, 0.01846862257356514, 0.01262296693390161, 1.28819427853, 7.729165671180001, 0.0001267527194229547, 0.00008663317362650538) },
    @{ Row = 15; D = "FAPs"; Vals = @(2, 0.6666666666666666, 0.1188766666666667, 0.35663, 0.006863138759702679, 0.006863138759702677, 3, 1, 139.6948166666666, 419.0844499999999, 0.2380839126543345, 0.2440890227431923, 16.60645415594444, 149.4580874035, 0.001634002928999631, 0.001675216832806752) },
    @{ Row = 16; D = "Inflammatory-Mac"; Vals = @(2, 0.6666666666666666, 0.1188766666666667, 0.35663, 0.006863138759702679, 0.006863138759702677, 3, 1, 150.12088, 450.36264, 0.2558532043948076, 0.2623065033256284, 17.84586981146667, 160.6128283032, 0.001755956043876136, 0.0018002459298962) },
    @{ Row = 17; D = "MuSCs"; Vals = @(2, 0.6666666666666666, 0.1188766666666667, 0.35663, 0.006863138759702679, 0.006863138759702677, 2, 1, 32.469223, 64.938446, 0.05533777012737728, 0.03782235733776705, 3.859832999496667, 23.15899799698, 0.00037979079503672, 0.000259580086628154) },
    @{ Row = 18; D = "Neutrophils"; Vals = @(2, 0.6666666666666666, 0.1188766666666667, 0.35663, 0.006863138759702679, 0.006863138759702677, 3, 1, 52.09024033333333, 156.270721, 0.08877815602319267, 0.09101737745760805, 6.192314136692222, 55.73082723022999, 0.0006092968036177054, 0.0006246648910357985) },
    @{ Row = 19; D = $null; Vals = @(2, 0.6666666666666666, 0.1188766666666667, 0.35663, 0.006863138759702679, 0.006863138759702677, 3, 1, 201.5345866666667, 604.60376, 0.3434783342267227, 0.3521417722019025, 23.95775988097778, 215.6198389288, 0.002357339468749532, 0.002416797845709268) }
)

foreach ($item in $rows) {
    $r = $item.Row
    if ($item.D -ne $null) {
        $ws.Range("D$r").Value = $item.D
    }
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $item.Vals[$i]
    }
}